$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Ensure we have 19 rows total (insert 4 blank rows before current row 16) ---
$ws.Rows.Item(16).Resize(4).Insert()

# --- Column B width change 51 -> 52 ---
# NOTE: Excels ColumnWidth (character units) and the OOXML <col width> differ
# by the fixed padding offset (5/6 for the default Calibri 11 font); subtract it
# here so the saved OOXML width attribute comes out to exactly 52.
$ws.Columns.Item(2).ColumnWidth = (52 - 5/6)

# --- Clear all existing hyperlinks (targets/order will be rebuilt below) ---
$ws.Hyperlinks.Delete()

# --- Write row data (A:H) ---
# row 2
$ws.Cells.Item(2,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(2,2).Value = '【SES案件多数】バックエンドエンジニア募集(Java/PHP/Python/Node.js)'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5399874'
$ws.Cells.Item(2,7).Value = 320
$ws.Cells.Item(2,8).Value = '🔥Python ★Java ◆Node.js ○PHP'
# row 3
$ws.Cells.Item(3,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(3,2).Value = '【低予算希望】LINE公式アカウント+社食注文システム開発依頼(社内利用のみ)'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5400375'
$ws.Cells.Item(3,7).Value = 118
$ws.Cells.Item(3,8).Value = '◆開発,システム開発'
# row 4
$ws.Cells.Item(4,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(4,2).Value = '【急募】LLMによるMCP(Model Context Protocol)でのExcel操作機能開発'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5400689'
$ws.Cells.Item(4,7).Value = 75
$ws.Cells.Item(4,8).Value = '◆開発'
# row 5
$ws.Cells.Item(5,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(5,2).Value = '【急募】住宅展示場マッチング診断サービスのMVP開発依頼'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5399759'
$ws.Cells.Item(5,7).Value = 75
$ws.Cells.Item(5,8).Value = '◆開発'
# row 6
$ws.Cells.Item(6,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(6,2).Value = '【フリーランス募集】CTビューアーソフト気道抽出機能開発'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5400101'
$ws.Cells.Item(6,7).Value = 68
$ws.Cells.Item(6,8).Value = '◆開発'
# row 7
$ws.Cells.Item(7,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(7,2).Value = '【急募】音源ライセンス販売サイトのMVP構築依頼'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5400763'
$ws.Cells.Item(7,7).Value = 45
$ws.Cells.Item(7,8).Value = '◇サイト'
# row 8
$ws.Cells.Item(8,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(8,2).Value = 'MYSQLからGoogleスプレッドシートへデータ取り込み及びスプレッドシート改修'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5400606'
$ws.Cells.Item(8,7).Value = 30
$ws.Cells.Item(8,8).Value = '◇MySQL'
# row 9
$ws.Cells.Item(9,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(9,2).Value = 'eBayテラピークでのキーワード検索結果等の取得するためのシステム制作'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5390238'
$ws.Cells.Item(9,7).Value = 33
$ws.Cells.Item(9,8).Value = ''
# row 10
$ws.Cells.Item(10,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(10,2).Value = 'Drupal関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5400683'
$ws.Cells.Item(10,7).Value = 25
$ws.Cells.Item(10,8).Value = ''
# row 11
$ws.Cells.Item(11,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(11,2).Value = '金融関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5400681'
$ws.Cells.Item(11,7).Value = 25
$ws.Cells.Item(11,8).Value = ''
# row 12
$ws.Cells.Item(12,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(12,2).Value = '【急募】SOLIDWORKS2024での機械設計と製図依頼'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5400338'
$ws.Cells.Item(12,7).Value = 25
$ws.Cells.Item(12,8).Value = ''
# row 13
$ws.Cells.Item(13,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(13,2).Value = '当社CTソフトへの機能追加:気道抽出'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '3,000,000 円 ~ 5,000,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5400094'
$ws.Cells.Item(13,7).Value = 25
$ws.Cells.Item(13,8).Value = ''
# row 14
$ws.Cells.Item(14,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(14,2).Value = '【SES案件多数/リモート可】フルスタックエンジニア募集(フロント〜バック〜クラウドまで歓迎)'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5399877'
$ws.Cells.Item(14,7).Value = 25
$ws.Cells.Item(14,8).Value = ''
# row 15
$ws.Cells.Item(15,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(15,2).Value = '【SES案件多数/リモート可】インフラエンジニア募集(AWS/Linux/NW設計・構築 等歓迎)'
$ws.Cells.Item(15,3).Value = 'システム開発'
$ws.Cells.Item(15,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(15,5).Value = '期限情報なし'
$ws.Cells.Item(15,6).Value = 'https://www.lancers.jp/work/detail/5399876'
$ws.Cells.Item(15,7).Value = 25
$ws.Cells.Item(15,8).Value = ''
# row 16
$ws.Cells.Item(16,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(16,2).Value = '【急募】Nuxt3でのWEBページ表示速度改善依頼'
$ws.Cells.Item(16,3).Value = 'システム開発'
$ws.Cells.Item(16,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(16,5).Value = '期限情報なし'
$ws.Cells.Item(16,6).Value = 'https://www.lancers.jp/work/detail/5400231'
$ws.Cells.Item(16,7).Value = 18
$ws.Cells.Item(16,8).Value = ''
# row 17
$ws.Cells.Item(17,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(17,2).Value = '限定公開 PR 限定公開の仕事'
$ws.Cells.Item(17,3).Value = 'システム開発'
$ws.Cells.Item(17,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(17,5).Value = '期限情報なし'
$ws.Cells.Item(17,6).Value = 'https://www.lancers.jp/work/detail/5399347'
$ws.Cells.Item(17,7).Value = 13
$ws.Cells.Item(17,8).Value = ''
# row 18
$ws.Cells.Item(18,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(18,2).Value = '【急募】ライフプランシミュレーターのバグ確認と使用感調査'
$ws.Cells.Item(18,3).Value = 'システム開発'
$ws.Cells.Item(18,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(18,5).Value = '期限情報なし'
$ws.Cells.Item(18,6).Value = 'https://www.lancers.jp/work/detail/5400626'
$ws.Cells.Item(18,7).Value = 10
$ws.Cells.Item(18,8).Value = ''
# row 19
$ws.Cells.Item(19,1).Value = '2025-09-26 01:15:09'
$ws.Cells.Item(19,2).Value = '【SalesIQ活用】CRMと連携したリード獲得方法を教えてください'
$ws.Cells.Item(19,3).Value = 'システム開発'
$ws.Cells.Item(19,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(19,5).Value = '期限情報なし'
$ws.Cells.Item(19,6).Value = 'https://www.lancers.jp/work/detail/5400402'
$ws.Cells.Item(19,7).Value = 10
$ws.Cells.Item(19,8).Value = ''

# --- Rebuild hyperlinks for column F (rows 2-19), in row order, for sequential rIds ---
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), 'https://www.lancers.jp/work/detail/5399874', "", "", 'https://www.lancers.jp/work/detail/5399874') | Out-Null
$ws.Cells.Item(2,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), 'https://www.lancers.jp/work/detail/5400375', "", "", 'https://www.lancers.jp/work/detail/5400375') | Out-Null
$ws.Cells.Item(3,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), 'https://www.lancers.jp/work/detail/5400689', "", "", 'https://www.lancers.jp/work/detail/5400689') | Out-Null
$ws.Cells.Item(4,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), 'https://www.lancers.jp/work/detail/5399759', "", "", 'https://www.lancers.jp/work/detail/5399759') | Out-Null
$ws.Cells.Item(5,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), 'https://www.lancers.jp/work/detail/5400101', "", "", 'https://www.lancers.jp/work/detail/5400101') | Out-Null
$ws.Cells.Item(6,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), 'https://www.lancers.jp/work/detail/5400763', "", "", 'https://www.lancers.jp/work/detail/5400763') | Out-Null
$ws.Cells.Item(7,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), 'https://www.lancers.jp/work/detail/5400606', "", "", 'https://www.lancers.jp/work/detail/5400606') | Out-Null
$ws.Cells.Item(8,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), 'https://www.lancers.jp/work/detail/5390238', "", "", 'https://www.lancers.jp/work/detail/5390238') | Out-Null
$ws.Cells.Item(9,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), 'https://www.lancers.jp/work/detail/5400683', "", "", 'https://www.lancers.jp/work/detail/5400683') | Out-Null
$ws.Cells.Item(10,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), 'https://www.lancers.jp/work/detail/5400681', "", "", 'https://www.lancers.jp/work/detail/5400681') | Out-Null
$ws.Cells.Item(11,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), 'https://www.lancers.jp/work/detail/5400338', "", "", 'https://www.lancers.jp/work/detail/5400338') | Out-Null
$ws.Cells.Item(12,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(13,6), 'https://www.lancers.jp/work/detail/5400094', "", "", 'https://www.lancers.jp/work/detail/5400094') | Out-Null
$ws.Cells.Item(13,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(14,6), 'https://www.lancers.jp/work/detail/5399877', "", "", 'https://www.lancers.jp/work/detail/5399877') | Out-Null
$ws.Cells.Item(14,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(15,6), 'https://www.lancers.jp/work/detail/5399876', "", "", 'https://www.lancers.jp/work/detail/5399876') | Out-Null
$ws.Cells.Item(15,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(16,6), 'https://www.lancers.jp/work/detail/5400231', "", "", 'https://www.lancers.jp/work/detail/5400231') | Out-Null
$ws.Cells.Item(16,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(17,6), 'https://www.lancers.jp/work/detail/5399347', "", "", 'https://www.lancers.jp/work/detail/5399347') | Out-Null
$ws.Cells.Item(17,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(18,6), 'https://www.lancers.jp/work/detail/5400626', "", "", 'https://www.lancers.jp/work/detail/5400626') | Out-Null
$ws.Cells.Item(18,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(19,6), 'https://www.lancers.jp/work/detail/5400402', "", "", 'https://www.lancers.jp/work/detail/5400402') | Out-Null
$ws.Cells.Item(19,6).Style = "Hyperlink"
